$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.637.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.97%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.27%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5155"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4388"
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.76%  "

# Row 11
$ws.Range("E11").Value = "  +0.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.099.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.32%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.759"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.161"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "100.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.29%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001156"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.38%  "

# Row 18
$ws.Range("E18").Value = "  -0.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.74%  "

# Row 20
$ws.Range("E20").Value = "  -0.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.180"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.684.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.37%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.316"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.337.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.60%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.98%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.530"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.04%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.84%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.13%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.134"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.57%  "

# Row 32
$ws.Range("E32").Value = "  -2.37%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.624"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.57%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.174"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.962"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.32%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.068"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.61%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02567"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06718"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.74%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2241"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.01%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6877"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.23%  "

# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6701"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.12%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.44%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.309"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.610"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.94%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.217"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000340"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.06%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.166"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.30%  "
